$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

$ws.Range("E1").Value = "schema.properties.latitude.type"
$ws.Range("F1").Value = "schema.properties.longitude.type"
$ws.Range("G1").Value = "schema.properties.altitude.type"
$ws.Range("H1").Value = "schema.properties.accuracy.type"
